$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated DM_Stat (C) and P_Value (D) columns for rows 2-11
$ws.Range("C2").Value = -1.291906760184714
$ws.Range("D2").Value = 0.2097967734461148

$ws.Range("C3").Value = -0.3840945598667236
$ws.Range("D3").Value = 0.7045956771925179

$ws.Range("C4").Value = -0.2355416140551202
$ws.Range("D4").Value = 0.815967180459285

$ws.Range("C5").Value = 0.633835692680039
$ws.Range("D5").Value = 0.5327233818121218

$ws.Range("C6").Value = 0.9163325101137673
$ws.Range("D6").Value = 0.3694285275095535

$ws.Range("C7").Value = 1.019615285723138
$ws.Range("D7").Value = 0.3189932889144209

$ws.Range("C8").Value = 2.02940847891387
$ws.Range("D8").Value = 0.05468546139969566

$ws.Range("C9").Value = 0.1029233583603398
$ws.Range("D9").Value = 0.918956230535966

$ws.Range("C10").Value = 0.7138602539101656
$ws.Range("D10").Value = 0.4828190257826499

$ws.Range("C11").Value = 0.594423221530923
$ws.Range("D11").Value = 0.5582932427371594

# With the updated P_Value, row 8 is no longer statistically significant
$ws.Range("G8").Value = "No"
